# Commit: "remove personal information before submit"
# Replace placeholder/internal codenames with the final public labels
# across both worksheets, and drop the leftover personal-info row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("MNC Experiment")
$ws2 = $wb.Worksheets.Item("CAH Experiment")

# ---- MNC Experiment sheet ----
$ws1.Range("C6").Value = "EI"
$ws1.Range("F7").Value = "LM2"
$ws1.Range("G7").Value = "LM3"
$ws1.Range("H7").Value = "LM4"

# ---- CAH Experiment sheet ----
$ws2.Range("I7").Value = "LM2"
$ws2.Range("J7").Value = "LM3"
$ws2.Range("K7").Value = "LM4"
$ws2.Range("G8").Value = "EI"
$ws2.Range("G24").Value = "Diff-GP"

# column D on CAH Experiment widened (was a narrow spacer column)
$ws2.Columns.Item(4).ColumnWidth = 11.33

# remove the leftover personal/internal label in E33
$ws2.Range("E33").ClearContents()

# restore selections to match the saved cursor position in each sheet
$ws1.Activate()
$ws1.Range("E32").Select()

$ws2.Activate()
$ws2.Range("G9").Select()
